$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-09-03 01:09:09"
$wsDeDe.Range("H2").Value = "2016-09-03 01:09:09"

$wsZhCn.Range("H2").Value = "2016-09-03 01:08:59"
$wsZhCn.Range("K2").Value = "2016-09-03 01:09:54"

$wsDeDe.Range("K2").Value = "2016-09-03 01:10:00"
